# Highlight two bullet items in green on the "Step 1" slide (slide 4):
#   - "Integrate code from Project 1 into Project 2" (paragraph 1)
#   - "Add and commit the files to version control" (paragraph 6)
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Green highlight color (00FF00) expressed as a VBA-style RGB long: R + G*256 + B*65536
$green = 65280

$tr.Paragraphs(1, 1).Font.Highlight.RGB = $green
$tr.Paragraphs(6, 1).Font.Highlight.RGB = $green
